$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 215210.6
$ws.Range("J17").Value = 215210.6
$ws.Range("L17").Value = 645631.8
$ws.Range("N17").Value = -645967.8

$ws.Range("H57").Value = 224962.17
$ws.Range("J57").Value = 224962.17
$ws.Range("L57").Value = 674886.51
$ws.Range("N57").Value = -675884.51

$ws.Range("H80").Value = 1009.2069
$ws.Range("I80").Value = 926.4167
$ws.Range("J80").Value = 1067.6471
$ws.Range("K80").Value = 2779.2501
$ws.Range("L80").Value = 3202.9413
$ws.Range("M80").Value = -1781.2501
$ws.Range("N80").Value = -5198.9413

$ws.Range("H83").Value = 1009.2069
$ws.Range("I83").Value = 926.4167
$ws.Range("J83").Value = 1067.6471
$ws.Range("K83").Value = 8337.7503
$ws.Range("L83").Value = 9608.823899999999
$ws.Range("M83").Value = -3345.7503
$ws.Range("N83").Value = -19592.8239

$ws.Range("H98").Value = 3644.5715
$ws.Range("I98").Value = 3644.5715
$ws.Range("K98").Value = 3644.5715
$ws.Range("M98").Value = -2146.5715

$ws.Range("H100").Value = 10441.235
$ws.Range("I100").Value = 1674.25
$ws.Range("K100").Value = 1674.25
$ws.Range("M100").Value = -1133.25

$ws.Range("H122").Value = 3644.5715
$ws.Range("I122").Value = 3644.5715
$ws.Range("K122").Value = 10933.7145
$ws.Range("M122").Value = -8483.7145

$ws.Range("H123").Value = 99000
$ws.Range("J123").Value = 99000
$ws.Range("L123").Value = 99000
$ws.Range("N123").Value = -108800

$ws.Range("H125").Value = 2765.8572
$ws.Range("I125").Value = 1465.25
$ws.Range("J125").Value = 4500
$ws.Range("K125").Value = 13187.25
$ws.Range("L125").Value = 40500
$ws.Range("M125").Value = -10727.25
$ws.Range("N125").Value = -45420

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3646.7083
$ws.Range("I45").Value = 3100.625
$ws.Range("K45").Value = 3100.625
$ws.Range("M45").Value = -2723.625

$ws.Range("H97").Value = 672.05884
$ws.Range("I97").Value = 719.6
$ws.Range("K97").Value = 719.6
$ws.Range("M97").Value = -223.6

$ws.Range("H122").Value = 3675.2766
$ws.Range("I122").Value = 3225.5
$ws.Range("K122").Value = 9676.5
$ws.Range("M122").Value = -7226.5

$ws.Range("H132").Value = 2217.4324
$ws.Range("I132").Value = 2120.4856
$ws.Range("K132").Value = 6361.4568
$ws.Range("M132").Value = -3831.4568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 892.5
$ws.Range("J80").Value = 925
$ws.Range("L80").Value = 925
$ws.Range("N80").Value = -2921

$ws.Range("H83").Value = 892.5
$ws.Range("J83").Value = 925
$ws.Range("L83").Value = 4625
$ws.Range("N83").Value = -14609

$ws.Range("H99").Value = 3301.5789
$ws.Range("I99").Value = 2203.923
$ws.Range("K99").Value = 2203.923
$ws.Range("M99").Value = -705.9229999999998

$ws.Range("H134").Value = 7405.706
$ws.Range("I134").Value = 3856.75
$ws.Range("K134").Value = 11570.25
$ws.Range("M134").Value = -9035.25

$ws.Range("H137").Value = 56250
$ws.Range("J137").Value = 56250
$ws.Range("L137").Value = 56250
$ws.Range("N137").Value = -66450

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 80001
$ws.Range("J64").Value = 80001
$ws.Range("L64").Value = 80001
$ws.Range("N64").Value = -80497

$ws.Range("H67").Value = 80001
$ws.Range("J67").Value = 80001
$ws.Range("L67").Value = 80001
$ws.Range("N67").Value = -81717

$ws.Range("H99").Value = 2222.0667
$ws.Range("I99").Value = 1945.909
$ws.Range("K99").Value = 1945.909
$ws.Range("M99").Value = -447.9090000000001

$ws.Range("H126").Value = 2222.0667
$ws.Range("I126").Value = 1945.909
$ws.Range("K126").Value = 5837.727000000001
$ws.Range("M126").Value = -3367.727000000001

$ws.Range("H134").Value = 4620
$ws.Range("I134").Value = 1853.8889
$ws.Range("K134").Value = 5561.6667
$ws.Range("M134").Value = -3026.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 254851
$ws.Range("J55").Value = 504000
$ws.Range("L55").Value = 1512000
$ws.Range("N55").Value = -1512354

$ws.Range("H57").Value = 4499.5
$ws.Range("I57").Value = 1500
$ws.Range("J57").Value = 7499
$ws.Range("K57").Value = 4500
$ws.Range("L57").Value = 22497
$ws.Range("M57").Value = -3941
$ws.Range("N57").Value = -23615

$ws.Range("H62").Value = 8256
$ws.Range("I62").Value = 7008
$ws.Range("J62").Value = 12000
$ws.Range("K62").Value = 21024
$ws.Range("L62").Value = 36000
$ws.Range("M62").Value = -20338
$ws.Range("N62").Value = -37372

$ws.Range("H65").Value = 8256
$ws.Range("I65").Value = 7008
$ws.Range("J65").Value = 12000
$ws.Range("K65").Value = 63072
$ws.Range("L65").Value = 108000
$ws.Range("M65").Value = -59640
$ws.Range("N65").Value = -114864

$ws.Range("H68").Value = 61980.43
$ws.Range("I68").Value = 107017
$ws.Range("K68").Value = 321051
$ws.Range("M68").Value = -320240

$ws.Range("H71").Value = 61980.43
$ws.Range("I71").Value = 107017
$ws.Range("K71").Value = 963153
$ws.Range("M71").Value = -959097

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 502
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

$ws.Range("H8").Value = 502
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()

$ws.Range("H53").Value = 39232
$ws.Range("J53").Value = 39232
$ws.Range("L53").Value = 39232
$ws.Range("N53").Value = -40494

$ws.Range("H97").Value = 1107.8529
$ws.Range("I97").Value = 596.7727
$ws.Range("K97").Value = 596.7727
$ws.Range("M97").Value = -100.7727

$ws.Range("H122").Value = 4193.1665
$ws.Range("I122").Value = 3172.9092
$ws.Range("J122").Value = 5056.4614
$ws.Range("K122").Value = 9518.7276
$ws.Range("L122").Value = 15169.3842
$ws.Range("M122").Value = -7068.7276
$ws.Range("N122").Value = -20069.3842

$ws.Range("H132").Value = 2854.5574
$ws.Range("I132").Value = 2646.5112
$ws.Range("K132").Value = 7939.5336
$ws.Range("M132").Value = -5409.5336

$ws.Range("H134").Value = 87494.25
$ws.Range("J134").Value = 87494.25
$ws.Range("L134").Value = 262482.75
$ws.Range("N134").Value = -267552.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3135.7273
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws.Range("H27").Value = 3135.7273
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()

$ws.Range("H93").Value = 2192.9375
$ws.Range("I93").Value = 2180.5386
$ws.Range("K93").Value = 2180.5386
$ws.Range("M93").Value = -932.5385999999999

$ws.Range("H122").Value = 3585.7144
$ws.Range("I122").Value = 1567
$ws.Range("J122").Value = 5099.75
$ws.Range("K122").Value = 4701
$ws.Range("L122").Value = 15299.25
$ws.Range("M122").Value = -2251
$ws.Range("N122").Value = -20199.25

$ws.Range("H132").Value = 3169.7222
$ws.Range("I132").Value = 3164.814
$ws.Range("K132").Value = 9494.441999999999
$ws.Range("M132").Value = -6964.441999999999

$ws.Range("H135").Value = 100000
$ws.Range("J135").Value = 100000
$ws.Range("L135").Value = 100000
$ws.Range("N135").Value = -110140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 30997
$ws.Range("I58").Value = 21995
$ws.Range("K58").Value = 21995
$ws.Range("M58").Value = -21687

$ws.Range("H62").Value = 71052.336
$ws.Range("J62").Value = 3999.5
$ws.Range("L62").Value = 3999.5
$ws.Range("N62").Value = -5247.5

$ws.Range("H65").Value = 71052.336
$ws.Range("J65").Value = 3999.5
$ws.Range("L65").Value = 19997.5
$ws.Range("N65").Value = -26237.5

$ws.Range("H100").Value = 1672.4445
$ws.Range("I100").Value = 1756.5
$ws.Range("K100").Value = 3513
$ws.Range("M100").Value = -2972
